$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: swap "Current_LT_Rating"/"B2" criteria pair for
#     "Corporate_Accounting_Regime"/"US GAAP" ---
$ws.Range("A9").Formula = '=DSLOOKUP("Monthly_Data_Query","Company","Alon USA Energy, Inc.","Year","2012","Corporate_Accounting_Regime","US GAAP","Country")'

# --- Row 10: same pair swap, different argument order ---
$ws.Range("A10").Formula = '=DSLOOKUP("Monthly_Data_Query","Company","Alon USA Energy, Inc.","Corporate_Accounting_Regime","US GAAP","Year","2012","Country")'

# --- Row 11: the criteria cells referenced by the formula move too ---
$ws.Range("F11").Value = "Corporate_Accounting_Regime"
$ws.Range("G11").Value = "US GAAP"
$ws.Range("H11").Value = "Year"
# I11 switches from the number 2012 to the text "2012" (quote-prefixed)
$ws.Range("I11").Value = "'2012"
$ws.Range("J11").Value = "Country"

# --- Row 13: different DB/criteria-field pair ---
$ws.Range("A13").Formula = '=DSLOOKUP("Monthly_Data_Query","Secured_Debt", 0, "GROSS_PROFIT", 4049115.256714, "Treasury_Stock_Issued_Repurchased")'

# --- Column A widened to fit the new, longer formula text ---
# (target stored width is 65.42578125; the engine quantizes ColumnWidth to
# 1/6-character steps, so 65.5 - the nearest reachable value - is what we get)
$ws.Columns("A").ColumnWidth = 64.65

# --- Selection moves to B17 ---
[void]$ws.Range("B17").Select()
